$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text updates derived from the cryptos-list refresh.
# Values that look numeric (e.g. "254.97") are force-written as text so they
# keep matching the source feed formatting instead of being auto-converted to
# Excel numbers; the cell Style is preserved (saved/restored) so no formatting
# side effects are introduced.

$ws.Range('D2').Value = '35.504.95'
$ws.Range('E2').Value = '  +0.69%  '
$ws.Range('D3').Value = '1.928.20'
$ws.Range('E3').Value = '  +1.71%  '
$ws.Range('E4').Value = '  -0.11%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.735'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  +11.87%  '
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '255.84'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  +5.04%  '
$ws.Range('E7').Value = '  -0.10%  '
$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '40.74'
$ws.Range('D8').Style = $origStyle
$ws.Range('E8').Value = '  -1.85%  '
$ws.Range('E9').Value = '  +3.74%  '
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.72'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  +4.71%  '
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0742'
$ws.Range('D11').Style = $origStyle
$ws.Range('E11').Value = '  +4.57%  '
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').Value = '2.206.39'
$ws.Range('E13').Value = '  +1.64%  '
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.81'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  +6.47%  '
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.719'
$ws.Range('D15').Style = $origStyle
$ws.Range('E15').Value = '  +3.58%  '
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$origStyle = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.92'
$ws.Range('D16').Style = $origStyle
$ws.Range('E16').Value = '  +1.97%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '1.913.76'
$ws.Range('E17').Value = '  +1.81%  '
$ws.Range('D18').Value = '35.474.81'
$ws.Range('E18').Value = '  +0.59%  '
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.97'
$ws.Range('D19').Style = $origStyle
$ws.Range('E19').Value = '  +3.82%  '
$ws.Range('D20').Value = '0.0₃0837'
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '242.53'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  +0.64%  '
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.05'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  +4.80%  '
$ws.Range('E23').Value = '  +8.25%  '
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('E25').Value = '  +2.14%  '
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.37'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  -1.16%  '
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '168.27'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  -1.15%  '
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.66'
$ws.Range('D28').Style = $origStyle
$ws.Range('E28').Value = '  +3.83%  '
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.135'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  +6.59%  '
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '18.88'
$ws.Range('D30').Style = $origStyle
$ws.Range('E30').Value = '  +3.70%  '
$ws.Range('D31').Value = '4.125.82'
$ws.Range('E31').Value = '  +19.37%  '
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.36'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  +6.17%  '
$ws.Range('E33').Value = '  +14.12%  '
$ws.Range('E34').Value = '  +23.35%  '
$ws.Range('E35').Value = '  +3.83%  '
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.27'
$ws.Range('D36').Style = $origStyle
$ws.Range('E36').Value = '  +4.14%  '
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('E38').Value = '  -2.45%  '
$ws.Range('E39').Value = '  +0.55%  '
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.26'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  +8.30%  '
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.27'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  +10.35%  '
$ws.Range('E42').Value = '  +4.72%  '
$ws.Range('E43').Value = '  +1.12%  '
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0650'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  +2.56%  '
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.48'
$ws.Range('D45').Style = $origStyle
$ws.Range('E45').Value = '  +4.92%  '
$ws.Range('D46').Value = '1.349.53'
$ws.Range('E46').Value = '  +0.99%  '
$ws.Range('E47').Value = '  +0.92%  '
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.78'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  +0.15%  '
$ws.Range('E49').Value = '  +3.31%  '
$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '45.39'
$ws.Range('D50').Style = $origStyle
$ws.Range('E50').Value = '  -4.16%  '
$ws.Range('B51').Value = 'Gas'
$ws.Range('C51').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.93'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  +3.13%  '

Write-Output "Applied 92 cell updates"
